# Finished the BOM for communication and sensors. Added the last obsticle
# detection sensor, an ESP32-CAM which includes a camera as well as a
# development board.

$wb = $excel.ActiveWorkbook

$obs = $wb.Worksheets.Item("Obstical detection BOM")
$complete = $wb.Worksheets.Item("Complete BOM")

# --- Move the existing "OPT8241NBN" secondary-component block from column D
# --- to column F (it becomes the tertiary component), then populate column D
# --- with the new ESP32-CAM primary replacement data.

# Column F (tertiary component) <- old column D content
$obs.Range("F1").Value = "OPT8241NBN"
$obs.Range("F2").Formula = "=HYPERLINK(""https://www.mouser.se/ProductDetail/Texas-Instruments/OPT8241NBN?qs=cGEy3R83DS%2FxFMUAL%252BoBvw%3D%3D"",""OPT8241NBN Texas Instruments | Mouser Sverige"")"
$obs.Range("F3").Value = "7,9*8,8*0,8"
$obs.Range("F4").Value = 1
$obs.Range("F5").Value = 605.13
$obs.Range("F6").Formula = "=F4*F5"

# Column D (secondary component) <- new ESP32-CAM data
$obs.Range("D1").Value = "ESP32-CAM Utvecklingskort"
$obs.Range("D2").Formula = "=HYPERLINK(""https://www.electrokit.com/en/esp32-cam-development-board"",""Electro:kit - ESP32 url"")"
$obs.Range("D3").Value = "(40*27*15)"
$obs.Range("D4").Value = 1
$obs.Range("D5").Value = 199
$obs.Range("D6").Formula = "=D4*D5"

# --- Column widths / cosmetic resizing that came along with the edit.
$obs.Columns.Item(2).ColumnWidth = 58.5546875
$obs.Columns.Item(6).ColumnWidth = 42.44140625

# --- Cursor/selection bookkeeping (cosmetic, matches the saved view state).
$obs.Range("D7").Select()
$complete.Range("G14").Select()

# --- The "Secondary Component" label in the Complete BOM roll-up sheet is
# --- cleared out as part of this pass.
$complete.Range("A9").Value = ""

$wb.Save()
